$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# All data cells in this sheet store plain text (inline strings) -
# coin names, coinranking.com links, prices (dot-grouped, e.g. "28.468.97")
# and hourly deltas ("  +0.19%  "). For price cells whose new text would
# otherwise be auto-converted to a number by Excel (e.g. "317.00" -> 317),
# we briefly force a Text number format so the literal string is kept,
# then restore the default "Normal" style so no extra formatting lingers.

$ws.Range("D2").Value = "28.468.97"
$ws.Range("E2").Value = "  +0.19%  "
$ws.Range("D3").Value = "1.795.56"
$ws.Range("E3").Value = "  -0.45%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.002"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  +0.18%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "317.00"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.29%  "
$ws.Range("E6").Value = "  +0.08%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.5403"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -0.49%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3781"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -1.11%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.07500"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -0.90%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "1.108"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -1.14%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "41.63"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -3.61%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.001"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -0.03%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "20.62"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -2.53%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.154"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -0.58%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "7.309"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +0.04%  "
$ws.Range("D16").Value = "1.805.35"
$ws.Range("E16").Value = "  +0.15%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "89.54"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -1.65%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.00001065"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -0.49%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.06486"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +0.35%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "17.44"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +0.94%  "
$ws.Range("E21").Value = "  +0.03%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.993"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +0.34%  "
$ws.Range("D23").Value = "28.492.85"
$ws.Range("E23").Value = "  +0.24%  "
$ws.Range("E24").Value = "  -0.77%  "
$ws.Range("E25").Value = "  -2.29%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "160.22"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +1.46%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "20.43"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -1.23%  "
$ws.Range("D28").Value = "1.997.68"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.310"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -4.90%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "122.55"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -1.22%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.108"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -4.07%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.1054"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +3.54%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "5.613"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -1.80%  "
$ws.Range("E34").Value = "  -0.18%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.2271"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +0.42%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.06528"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +4.59%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.02303"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -0.83%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "5.029"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +0.43%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "8.547"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -3.75%  "
$ws.Range("B40").Value = "TrustWalletToken"
$ws.Range("C40").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.200"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +3.91%  "
$ws.Range("B41").Value = "TheSandbox"
$ws.Range("C41").Value = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.6204"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -2.54%  "
$ws.Range("B42").Value = "WEMIXTOKEN"
$ws.Range("C42").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.450"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +5.23%  "
$ws.Range("B43").Value = "Aptos"
$ws.Range("C43").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "11.15"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -3.31%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "1.000"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -0.07%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "13.38"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -0.13%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "3.684"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +0.27%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.5810"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -2.40%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "126.85"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +2.63%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.948"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -0.46%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.192"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +4.02%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.06877"
$ws.Range("D51").Style = "Normal"
